$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Range("A20").Value = 111927215
$ws.Range("I20").Value = '''10'
$ws.Range("J20").Value = 'fruktkroppar'
$ws.Range("Q20").Value = 663485.6413922446
$ws.Range("R20").Value = 6602647.390513759

# Row 22
$ws.Range("A22").Value = 111926622
$ws.Range("I22").Value = ""
$ws.Range("J22").Value = ""
$ws.Range("Q22").Value = 663452.3464515609
$ws.Range("R22").Value = 6602675.90838708

# Row 23
$ws.Range("A23").Value = 112084535
$ws.Range("B23").Value = 88899
$ws.Range("D23").Value = 'NT'
$ws.Range("E23").Value = 3286
$ws.Range("F23").Value = 'Flattoppad klubbsvamp'
$ws.Range("G23").Value = 'Clavariadelphus truncatus'
$ws.Range("H23").Value = '(Quél.) Donk'
$ws.Range("I23").Value = '''80'
$ws.Range("J23").Value = 'fruktkroppar'
$ws.Range("Q23").Value = 663374.2695844367
$ws.Range("R23").Value = 6602611.054278261
$ws.Range("S23").Value = 10
$ws.Range("Z23").Value = '11:46'
$ws.Range("AB23").Value = '11:46'
$ws.Range("AC23").Value = 'Sötaktig mild smak (ej bitter)'

# Row 24
$ws.Range("A24").Value = 112084040
$ws.Range("I24").Value = ""
$ws.Range("J24").Value = ""
$ws.Range("Q24").Value = 663584.9559231531
$ws.Range("R24").Value = 6602703.52117154

# Row 25
$ws.Range("A25").Value = 112084114
$ws.Range("I25").Value = '''20'
$ws.Range("J25").Value = 'plantor/tuvor'
$ws.Range("K25").Value = 'fullt utvecklade blad'
$ws.Range("Q25").Value = 663576.8087203993
$ws.Range("R25").Value = 6602715.356141716
$ws.Range("S25").Value = 5

# Row 27
$ws.Range("A27").Value = 112083991
$ws.Range("B27").Value = 98535
$ws.Range("D27").Value = 'LC'
$ws.Range("E27").Value = 222498
$ws.Range("F27").Value = 'Blåsippa'
$ws.Range("G27").Value = 'Hepatica nobilis'
$ws.Range("H27").Value = 'Schreb.'
$ws.Range("I27").Value = '''300'
$ws.Range("J27").Value = 'stjälkar/strån/skott'
$ws.Range("Q27").Value = 663568.3519142884
$ws.Range("R27").Value = 6602664.1969273
$ws.Range("Z27").Value = '00:00'
$ws.Range("AB27").Value = '00:00'
$ws.Range("AC27").Value = ""

# Row 28
$ws.Range("A28").Value = 112083804
$ws.Range("I28").Value = '''300'
$ws.Range("J28").Value = 'stjälkar/strån/skott'
$ws.Range("K28").Value = ""
$ws.Range("Q28").Value = 663571.7306570449
$ws.Range("R28").Value = 6602738.498618284
$ws.Range("S28").Value = 10

# Row 29
$ws.Range("A29").Value = 112083737
$ws.Range("I29").Value = '''200'
$ws.Range("J29").Value = 'stjälkar/strån/skott'
$ws.Range("K29").Value = ""
$ws.Range("Q29").Value = 663545.1917381487
$ws.Range("R29").Value = 6602752.072187248
$ws.Range("S29").Value = 30

# Row 30
$ws.Range("A30").Value = 112083958
$ws.Range("I30").Value = '''10'
$ws.Range("J30").Value = 'plantor/tuvor'
$ws.Range("K30").Value = 'fullt utvecklade blad'
$ws.Range("Q30").Value = 663551.019940288
$ws.Range("R30").Value = 6602700.011799707
$ws.Range("S30").Value = 5
